# Update the "Generate Report for Handback" timestamps across the
# Overview, zh-cn and de-de sheets of the handback status workbook.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first file.
$wsOverview.Range("G2").Value = "2016-08-29 19:12:00"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# for the first file.
$wsZhCn.Range("H2").Value = "2016-08-29 19:11:55"
$wsZhCn.Range("K2").Value = "2016-08-29 19:12:31"

# de-de sheet: "Correspond Handoff Datetime" (mirrors the Overview value for
# the same file) and "Correspond Handback DateTime" for the first file.
$wsDeDe.Range("H2").Value = "2016-08-29 19:12:00"
$wsDeDe.Range("K2").Value = "2016-08-29 19:12:38"
